$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04) for every data
# row (rows 2 through 247).
$lastRow = 247
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}
